$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column keeps its original text representation by
# formatting the cells as Text before writing numeric-looking strings,
# so Excel does not coerce values like "134.67" into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.096.83"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.507.29"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "536.29"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "134.67"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").Value = "2.513.17"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "2.953.57"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "58.886.23"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "22.46"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "2.514.44"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "322.04"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "6.28"
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "65.86"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "7.44"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "174.49"
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("D30").Value = "0.0₃0761"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "6.28"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "18.13"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").Value = "3.94"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  +4.02%  "
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  +6.69%  "
$ws.Range("D41").Value = "36.61"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "276.48"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "131.77"
$ws.Range("E44").Value = "  +8.57%  "
$ws.Range("D45").Value = "5.06"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "0.0942"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("D48").Value = "0.0510"
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").Value = "0.0219"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").Value = "16.95"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "1.754.55"
$ws.Range("E51").Value = "  +0.78%  "
